$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J, copying the header style (bold, bordered) from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-72: column I ("I0") and column J ("IF") values
$data = @(
    @(2,7,8),
    @(3,6,6),
    @(4,5,5),
    @(5,5,5),
    @(6,5,6),
    @(7,7,8),
    @(8,6,6),
    @(9,6,6),
    @(10,9,9),
    @(11,6,7),
    @(12,10,10),
    @(13,7,8),
    @(14,7,7),
    @(15,9,9),
    @(16,7,8),
    @(17,7,7),
    @(18,8,8),
    @(19,9,9),
    @(20,6,7),
    @(21,7,7),
    @(22,5,6),
    @(23,7,8),
    @(24,5,6),
    @(25,7,7),
    @(26,3,4),
    @(27,6,7),
    @(28,8,8),
    @(29,7,8),
    @(30,7,7),
    @(31,5,6),
    @(32,8,8),
    @(33,8,8),
    @(34,8,9),
    @(35,7,8),
    @(36,8,8),
    @(37,8,8),
    @(38,7,7),
    @(39,6,6),
    @(40,7,7),
    @(41,8,8),
    @(42,6,7),
    @(43,7,7),
    @(44,7,7),
    @(45,8,8),
    @(46,7,7),
    @(47,8,8),
    @(48,7,7),
    @(49,7,7),
    @(50,7,7),
    @(51,7,7),
    @(52,6,7),
    @(53,8,8),
    @(54,5,5),
    @(55,6,6),
    @(56,5,6),
    @(57,6,7),
    @(58,6,6),
    @(59,7,7),
    @(60,6,6),
    @(61,6,6),
    @(62,7,7),
    @(63,5,6),
    @(64,7,7),
    @(65,8,8),
    @(66,8,8),
    @(67,2,5),
    @(68,9,9),
    @(69,1,3),
    @(70,1,3),
    @(71,4,5),
    @(72,6,6)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}

Write-Output "Added columns I (I0) and J (IF) for rows 1-72"
